$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '35.699.33'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.47%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.985.44'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -4.01%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.95'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.640'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.55%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '57.38'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +8.53%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '59.77'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0732'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('E12').Value = '  -4.58%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.927'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.55%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '14.11'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.77%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.276.33'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.93%  '
$ws.Range('E16').Value = '  -2.69%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.989.34'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.45%  '
$ws.Range('E18').Value = '  +5.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '35.568.03'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '70.63'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.55%  '
$ws.Range('E21').Value = '  -2.93%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '233.36'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('E23').Value = '  -4.07%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.35'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +10.39%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.15'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '163.65'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.52'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.77%  '
$ws.Range('E30').Value = '  -3.37%  '
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('E32').Value = '  -5.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0590'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0901'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +9.86%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.38'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.52%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.26'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -6.85%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.81'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.94'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('E40').Value = '  -4.93%  '
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0211'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.42%  '
$ws.Range('E43').Value = '  -4.22%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '91.14'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.10%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0888'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.79%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.380.02'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.45'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.54%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '15.49'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('E50').Value = '  -2.82%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '45.86'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.36%  '
